$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Replace "Jesus" with "Daniel" as the host for the two Monday office-hour slots
$ws.Range("B5").Value = "Daniel"
$ws.Range("B6").Value = "Daniel"

# Update the active selection on the sheet
$ws.Range("B7").Select()
